$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 ("STM32H743IIT6" / "Main MCU ALT" / "REMOVE", with an LCSC
# hyperlink on F2) was deleted from the BOM. Deleting the row shifts every
# row below it up by one.
$ws.Rows.Item(2).Delete()

# That row carried the only hyperlink in the sheet; remove it explicitly
# (row-delete alone leaves the link attached to whatever is now at F2).
$ws.Hyperlinks.Delete()

# A new part was added at the bottom of the table: a 26 MHz oscillator
# (KDS Daishinku). Fill Manufacturer before MFR. Part # so the shared
# string table is populated in the same order as the source edit.
$ws.Range("C13").Value = "KDS Daishinku"
$ws.Range("B13").Value = "1C/N226000AA0D"
$ws.Range("D13").Value = "RF"
$ws.Range("E13").Value = "SMD-3225_4P"
$ws.Range("F13").Value = "C160424"
$ws.Range("G13").Value = "26 MHz Oscillator"

# Match the "section divider" look used by the other first-row-of-a-group
# entries (e.g. row 11, the Storage section) by copying its formatting
# onto the new row's designator cell, then give the row the same height.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = 15

# Reset the view: no frozen/scrolled top-left cell, selection on A16.
$w = $excel.ActiveWindow
$w.ScrollRow = 1
$w.ScrollColumn = 1
$ws.Range("A16").Select()
